# Tuntikirjanpito.xlsx update
# - third day of work logged (row 22 gets hours+description, rows 23-30 added)
# - header "aika" renamed to "aika(h)"
# - totals row moved from 30 to 31, now summing B2:B30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New entries for the third day (2021-12-01, serial 44532) ---
# Row 22 already has the date in A22; add the hours + description.
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "päivämäärä-listan tarkistus ja korjaus, 00:00 muuttui automaattisesti muotoon 23:00"

$ws.Range("B23").Value = 1.5
$ws.Range("C23").Value = "ensimmäinen datan muokkausfunktio tehty, yksi datapiste per päivä, lähin keskiyöltä, getOneDataPointPerDate"

$ws.Range("B24").Value = 1
$ws.Range("C24").Value = "toinen datan muokkausfuntio tehty, pisin 'bearish trend' haetualla aikavälillä, getLongestBearishTrend"

$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "kolmas muokkausfuntio tehty, suurin volyymi päivä, funktioiden pientä refaktotorointia, forEach ==> map, tarkemmat nimet"

# --- Header update: "aika" -> "aika(h)" ---
$ws.Range("B1").Value = "aika(h)"

$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "neljäs muokkausfuntio tehty, paras päivä ostaa ja myydä, eli suurin profit"

$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "edellisen funktion korjausta testausta eri sarjoilla"

$ws.Range("B29").Value = 1
$ws.Range("C29").Value = "same date error lisätty, perus data näkyy bearish trend, highest trading volume, best buy/sell"

$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "Dataview sisällön perustan luonti, hardcoded data testaukseen"

# Row 30 used to be the "tunnit yht." total row; repurpose it as a normal data row.
$ws.Range("A30").Clear()
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = "refaktorointia uusien komponenttien kanssa, funktioiden testausta, kommentoinnin korjausta ja lisäystä"

# --- New total row 31 ---
$ws.Range("A31").Value = "tunnit yht."
$ws.Range("A31").HorizontalAlignment = -4108
$ws.Range("B31").Formula = "=SUM(B2:B30)"

# --- Cosmetic updates ---
$ws.Columns.Item(3).ColumnWidth = 114.3
$ws.Range("C30").Select() | Out-Null
